$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "CodeGPTPy"
$ws.Range("B9").Value = 32
$ws.Range("C9").Value = "Adam"
$ws.Range("D9").Value = 0

$ws.Range("E9").Select()
